$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Add the six new "electricity <X>g/kWh" fuel rows (25-30), mirroring the
#    existing "electricity PROXY" rows (B=1, C=1, D=0 with the D11/D12 style).
# ---------------------------------------------------------------------------
$names = @(
    "electricity 0g/kWh",
    "electricity 200g/kWh",
    "electricity 400g/kWh",
    "electricity 600g/kWh",
    "electricity 800g/kWh",
    "electricity 1000g/kWh"
)

$startRow = 25
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = 1
    $ws.Cells.Item($r, 3).Value = 1
    $ws.Cells.Item($r, 4).Value = 0
}

# Give the new D cells the same number-format style already used by D11/D12.
$null = $ws.Cells.Item(11, 4).Copy()
$null = $ws.Range("D25:D30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# A29/A30 pick up a distinct font (sz 11, black, Calibri) not used elsewhere
# in the sheet yet -- apply it once, then propagate via a format-only copy so
# only a single new font/style entry is minted.
$ws.Cells.Item(29, 1).Font.Size = 11
$ws.Cells.Item(29, 1).Font.Color = 0
$null = $ws.Cells.Item(29, 1).Copy()
$null = $ws.Range("A30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Re-enter F14:F23 as one shared formula (matches a fill-down over the
#    existing range).
# ---------------------------------------------------------------------------
$ws.Range("F14:F23").Formula2 = "=D14*`$F`$2"

# ---------------------------------------------------------------------------
# 3. Update the active selection to F28, like the saved workbook shows.
# ---------------------------------------------------------------------------
$null = $ws.Range("F28").Select()
